$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 233. This shifts the existing weekly blocks
# (previously at rows 233-276) down to rows 237-280, and leaves rows
# 233-236 free for the new "44476" weekly block.
$ws.Rows("233:236").Insert()

# Fill in the new weekly block (date 44476) in rows 233-236, reusing the
# same reference/location/price-per-kg metadata as the block that used to
# sit there (now shifted to rows 237-240), only Volumen (J) differs.
$newRows = @(
    @{ Row = 233; Calidad = "Especial"; Vol = 300; Min = 13500; Max = 14000; Prom = 13750; PKg = 764 },
    @{ Row = 234; Calidad = "Primera";  Vol = 500; Min = 11500; Max = 12000; Prom = 11750; PKg = 653 },
    @{ Row = 235; Calidad = "Segunda";  Vol = 360; Min = 9500;  Max = 10000; Prom = 9750;  PKg = 542 },
    @{ Row = 236; Calidad = "Tercera";  Vol = 240; Min = 6500;  Max = 7000;  Prom = 6750;  PKg = 375 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44476
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = 100112043
    $ws.Cells.Item($row, 7).Value = "Pepino dulce"
    $ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Vol
    $ws.Cells.Item($row, 11).Value = $r.Min
    $ws.Cells.Item($row, 12).Value = $r.Max
    $ws.Cells.Item($row, 13).Value = $r.Prom
    $ws.Cells.Item($row, 14).Value = "$/bandeja 18 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
